$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 2 4 '90.225.75'
Set-TextValue 2 5 '  -0.09%  '
Set-TextValue 3 4 '3.083.58'
Set-TextValue 3 5 '  -0.16%  '
Set-TextValue 4 5 '  +0.15%  '
Set-TextValue 5 4 '242.22'
Set-TextValue 5 5 '  +3.77%  '
Set-TextValue 6 4 '618.71'
Set-TextValue 6 5 '  -1.05%  '
Set-TextValue 7 5 '  +2.56%  '
Set-TextValue 8 4 '0.363'
Set-TextValue 8 5 '  +0.46%  '
Set-TextValue 9 5 '  +0.07%  '
Set-TextValue 10 4 '3.086.56'
Set-TextValue 10 5 '  +23.56%  '
Set-TextValue 11 4 '0.737'
Set-TextValue 11 5 '  +1.02%  '
Set-TextValue 12 5 '  +3.41%  '
Set-TextValue 13 4 '0.0000246'
Set-TextValue 13 5 '  -0.45%  '
Set-TextValue 14 4 '34.83'
Set-TextValue 14 5 '  -4.49%  '
Set-TextValue 15 5 '  -0.29%  '
Set-TextValue 16 4 '90.196.25'
Set-TextValue 17 4 '3.666.49'
Set-TextValue 17 5 '  +0.38%  '
Set-TextValue 18 4 '3.115.26'
Set-TextValue 18 5 '  +1.22%  '
Set-TextValue 19 4 '3.73'
Set-TextValue 19 5 '  -0.69%  '
Set-TextValue 20 5 '  +4.47%  '
Set-TextValue 21 4 '0.0000208'
Set-TextValue 21 5 '  -1.93%  '
Set-TextValue 22 4 '5.79'
Set-TextValue 22 5 '  +4.00%  '
Set-TextValue 23 4 '437.34'
Set-TextValue 23 5 '  -0.22%  '
Set-TextValue 24 5 '  +1.64%  '
Set-TextValue 25 5 '  -3.99%  '
Set-TextValue 26 4 '89.21'
Set-TextValue 26 5 '  +0.15%  '
Set-TextValue 27 4 '11.79'
Set-TextValue 27 5 '  -3.01%  '
Set-TextValue 28 4 '3.245.04'
Set-TextValue 29 4 '0.999'
Set-TextValue 29 5 '  -0.15%  '
Set-TextValue 30 5 '  +19.34%  '
Set-TextValue 31 4 '0.175'
Set-TextValue 31 5 '  +9.86%  '
Set-TextValue 32 4 '0.121'
Set-TextValue 32 5 '  +36.05%  '
Set-TextValue 33 4 '9.17'
Set-TextValue 33 5 '  -3.25%  '
Set-TextValue 34 5 '  +9.57%  '
Set-TextValue 35 2 'RenderToken'
Set-TextValue 35 3 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-TextValue 35 4 '7.92'
Set-TextValue 35 5 '  +12.99%  '
Set-TextValue 36 2 'MantraDAO'
Set-TextValue 36 3 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
Set-TextValue 36 4 '4.40'
Set-TextValue 36 5 '  +25.03%  '
Set-TextValue 37 2 'EthereumClassic'
Set-TextValue 37 3 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 37 4 '26.16'
Set-TextValue 37 5 '  -0.12%  '
Set-TextValue 38 2 'PancakeSwap'
Set-TextValue 38 3 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 38 4 '1.90'
Set-TextValue 38 5 '  -0.09%  '
Set-TextValue 39 2 'Bittensor'
Set-TextValue 39 3 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 39 4 '490.96'
Set-TextValue 39 5 '  -2.59%  '
Set-TextValue 40 2 'dogwifhat'
Set-TextValue 40 3 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue 40 4 '3.55'
Set-TextValue 40 5 '  -7.58%  '
Set-TextValue 41 2 'Fetch.AI'
Set-TextValue 41 3 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 41 4 '1.29'
Set-TextValue 41 5 '  +0.52%  '
Set-TextValue 42 2 'PolygonEcosystemToken'
Set-TextValue 42 3 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue 42 4 '0.419'
Set-TextValue 42 5 '  +2.03%  '
Set-TextValue 43 2 'WhiteBITCoin'
Set-TextValue 43 3 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue 43 4 '22.16'
Set-TextValue 43 5 '  -0.11%  '
Set-TextValue 44 2 'USDe'
Set-TextValue 44 3 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue 44 4 '1.00'
Set-TextValue 44 5 '  -0.02%  '
Set-TextValue 45 2 'Stacks'
Set-TextValue 45 3 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 45 4 '1.91'
Set-TextValue 45 5 '  +0.42%  '
Set-TextValue 46 2 'Monero'
Set-TextValue 46 3 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 46 4 '153.81'
Set-TextValue 46 5 '  +2.05%  '
Set-TextValue 47 2 'ARBITRUM'
Set-TextValue 47 3 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 47 4 '0.687'
Set-TextValue 47 5 '  -0.44%  '
Set-TextValue 48 2 'Binance-PegBSC-USD'
Set-TextValue 48 3 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue 48 4 '0.646'
Set-TextValue 48 5 '  -34.17%  '
Set-TextValue 49 5 '  -0.12%  '
Set-TextValue 51 4 '1.00'
Set-TextValue 51 5 '  +0.02%  '
